# The "trt" (treatment) column had its two group labels rewritten:
#   - what used to read "Placebo" (rows 2-10) is now the short label "PLA"
#   - what used to read "GTE" (rows 11-19) keeps reading "GTE"
# (net effect: the "Placebo" label is renamed to "PLA"; the cells that were
# already "GTE" are re-affirmed as "GTE")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B2:B10").Value = "PLA"
$ws.Range("B11:B19").Value = "GTE"

# author's cursor ended up on B10 when the file was saved
$ws.Range("B10").Select()
